$d = $word.ActiveDocument

$d.Content.Find.Execute("That's great, we will now start with the chatbot. ", $true, $false, $false, $false, $false, $true, 1, $false, "رائع،سنبدأ الآن مع روبوت الدردشة", 2)

$d.Content.Find.Execute("This chatbot will give you useful information about how to respond in a crisis. ", $true, $false, $false, $false, $false, $true, 1, $false, "سيزودك روبوت الدردشة بمعلومات مفيدة حول كيفية التعامل مع الأزمات.", 2)

$d.Content.Find.Execute("What's your child's name?", $true, $false, $false, $false, $false, $true, 1, $false, "ما اسم طفلك؟", 2)

$d.Content.Find.Execute("How old is @results.childname?", $true, $false, $false, $false, $false, $true, 1, $false, "ما عمر @results.childname؟ ", 2)

$d.Content.Find.Execute("Does @results.childname have siblings?", $true, $false, $false, $false, $false, $true, 1, $false, "هل  @results.childname عنده اخوة؟ ", 2)

$d.Content.Find.Execute("@results.childname likes football?", $true, $false, $false, $false, $false, $true, 1, $false, "@results.childname بحب الكرة ؟ ", 2)

$d.Content.Find.Execute("How old is {childname}?", $true, $false, $false, $false, $false, $true, 1, $false, "كم عمر {childname}؟", 2)

$d.Content.Find.Execute("Does {childname} have siblings?", $true, $false, $false, $false, $false, $true, 1, $false, "هل {childname} عنده اخوه ؟ ", 2)

$d.Content.Find.Execute("{childname} likes football?", $true, $false, $false, $false, $false, $true, 1, $false, "{childname} بحب الكرة ؟ ", 2)
